$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are stored as text in the source data (some contain
# thousands-separator dots like '26.339.83', others would lose trailing zeros or
# flip to scientific notation if Excel auto-typed them as numbers), so force the
# cell format to Text before writing each one.
$dCells = @("D2","D3","D5","D6","D8","D9","D11","D12","D13","D14","D15","D16","D17","D18","D19","D21","D22","D23","D24","D25","D26","D28","D29","D31","D32","D33","D34","D35","D36","D37","D39","D40","D41","D43","D44","D45","D46","D47","D48","D50","D51")
foreach ($cell in $dCells) { $ws.Range($cell).NumberFormat = "@" }

$ws.Range('D2').Value = '26.339.83'
$ws.Range('D3').Value = '1.690.43'
$ws.Range('D5').Value = '218.64'
$ws.Range('D6').Value = '0.5277'
$ws.Range('D8').Value = '0.2710'
$ws.Range('D9').Value = '22.13'
$ws.Range('D11').Value = '0.07503'
$ws.Range('D12').Value = '1.720.60'
$ws.Range('D13').Value = '4.567'
$ws.Range('D14').Value = '0.5851'
$ws.Range('D15').Value = '0.000008528'
$ws.Range('D16').Value = '64.59'
$ws.Range('D17').Value = '26.378.27'
$ws.Range('D18').Value = '4.946'
$ws.Range('D19').Value = '1.007'
$ws.Range('D21').Value = '189.64'
$ws.Range('D22').Value = '6.224'
$ws.Range('D23').Value = '1.007'
$ws.Range('D24').Value = '144.82'
$ws.Range('D25').Value = '7.716'
$ws.Range('D26').Value = '0.1240'
$ws.Range('D28').Value = '0.06665'
$ws.Range('D29').Value = '1.357'
$ws.Range('D31').Value = '3.589'
$ws.Range('D32').Value = '3.580'
$ws.Range('D33').Value = '1.673'
$ws.Range('D34').Value = '1.032'
$ws.Range('D35').Value = '0.6258'
$ws.Range('D36').Value = '2.393'
$ws.Range('D37').Value = '2.702'
$ws.Range('D39').Value = '1.118.18'
$ws.Range('D40').Value = '0.01624'
$ws.Range('D41').Value = '0.8907'
$ws.Range('D43').Value = '100.90'
$ws.Range('D44').Value = '1.839.21'
$ws.Range('D45').Value = '0.00000000116'
$ws.Range('D46').Value = '57.02'
$ws.Range('D47').Value = '8.195'
$ws.Range('D48').Value = '1.007'
$ws.Range('D50').Value = '6.108'
$ws.Range('D51').Value = '0.4302'

# Volume(1h) column (E) values are plain text percentages (e.g. '  +0.52%  ')
# and Excel never auto-converts these to numbers, so a direct assignment is safe.
$ws.Range('E2').Value = '  +0.52%  '
$ws.Range('E3').Value = '  +1.41%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('E5').Value = '  +0.40%  '
$ws.Range('E6').Value = '  +4.32%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  +2.25%  '
$ws.Range('E9').Value = '  +3.13%  '
$ws.Range('E10').Value = '  +1.28%  '
$ws.Range('E11').Value = '  +2.01%  '
$ws.Range('E12').Value = '  +3.12%  '
$ws.Range('E13').Value = '  +0.67%  '
$ws.Range('E14').Value = '  +0.91%  '
$ws.Range('E15').Value = '  +0.54%  '
$ws.Range('E16').Value = '  -0.01%  '
$ws.Range('E17').Value = '  +0.35%  '
$ws.Range('E18').Value = '  +0.59%  '
$ws.Range('E19').Value = '  +0.01%  '
$ws.Range('E20').Value = '  +1.00%  '
$ws.Range('E21').Value = '  +0.75%  '
$ws.Range('E22').Value = '  +0.79%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('E24').Value = '  +1.00%  '
$ws.Range('E25').Value = '  +0.74%  '
$ws.Range('E26').Value = '  +5.74%  '
$ws.Range('E27').Value = '  +1.51%  '
$ws.Range('E28').Value = '  +14.62%  '
$ws.Range('E29').Value = '  +6.15%  '
$ws.Range('E30').Value = '  +0.65%  '
$ws.Range('E31').Value = '  +2.35%  '
$ws.Range('E32').Value = '  +1.52%  '
$ws.Range('E33').Value = '  +2.54%  '
$ws.Range('E34').Value = '  +2.27%  '
$ws.Range('E35').Value = '  +4.62%  '
$ws.Range('E36').Value = '  +1.51%  '
$ws.Range('E37').Value = '  +2.30%  '
$ws.Range('E38').Value = '  +5.96%  '
$ws.Range('E39').Value = '  +4.46%  '
$ws.Range('E40').Value = '  +0.90%  '
$ws.Range('E41').Value = '  +3.10%  '
$ws.Range('E43').Value = '  +1.41%  '
$ws.Range('E44').Value = '  +1.13%  '
$ws.Range('E45').Value = '  +4.43%  '
$ws.Range('E46').Value = '  +2.38%  '
$ws.Range('E47').Value = '  +1.46%  '
$ws.Range('E48').Value = '  +0.28%  '
$ws.Range('E49').Value = '  +1.69%  '
$ws.Range('E50').Value = '  +4.40%  '
$ws.Range('E51').Value = '  +0.15%  '
